# Append two new job-listing rows to the "ランサーズ" sheet and refresh the
# "取得日時" (fetched-at) timestamp on every existing data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-02-11 02:06:47"

# Update column A (取得日時) for the existing data rows (2-18) to the new
# fetch timestamp.
for ($r = 2; $r -le 18; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# Row 19
$ws.Cells.Item(19, 1).Value = $newTimestamp
$ws.Cells.Item(19, 2).Value = "【急募】ドメイン接続業務の専門家を探しています!"
$ws.Cells.Item(19, 3).Value = "システム開発"
$ws.Cells.Item(19, 4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(19, 5).Value = "期限情報なし"
$ws.Cells.Item(19, 6).Value = "https://www.lancers.jp/work/detail/5489674"
$ws.Hyperlinks.Add($ws.Cells.Item(19, 6), "https://www.lancers.jp/work/detail/5489674")
$ws.Cells.Item(19, 6).Style = "Hyperlink"
$ws.Cells.Item(19, 7).Value = 18

# Row 20
$ws.Cells.Item(20, 1).Value = $newTimestamp
$ws.Cells.Item(20, 2).Value = "【急募】google work space の設定を専門家に依頼したい"
$ws.Cells.Item(20, 3).Value = "システム開発"
$ws.Cells.Item(20, 4).Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Cells.Item(20, 5).Value = "期限情報なし"
$ws.Cells.Item(20, 6).Value = "https://www.lancers.jp/work/detail/5489636"
$ws.Hyperlinks.Add($ws.Cells.Item(20, 6), "https://www.lancers.jp/work/detail/5489636")
$ws.Cells.Item(20, 6).Style = "Hyperlink"
$ws.Cells.Item(20, 7).Value = 10
